$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A - remaining columns (old B:F) shift left to A:E
$ws.Columns.Item(1).Delete()

# Rename header text "MODEL_CONDITION" -> "MODELCONDITION" (now in column D after the shift)
$ws.Cells.Item(1, 4).Value = "MODELCONDITION"
